$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove one of the two consecutive, identical, empty paragraphs that
#    sit right before the "github.com/conceptslearningmachine" line.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$targetIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("github.com/conceptslearningmachine")) {
        $targetIdx = $i
        break
    }
}

if ($targetIdx -gt 2) {
    $delIdx = $targetIdx - 2
    $checkIdx = $targetIdx - 1
    $delPara = $d.Paragraphs.Item($delIdx)
    $checkPara = $d.Paragraphs.Item($checkIdx)

    # Only the paragraph mark (no other text) should be present in both,
    # matching the duplicate blank-paragraph pattern being collapsed.
    if ($delPara.Range.Text.Length -eq 1 -and $checkPara.Range.Text.Length -eq 1) {
        $delPara.Range.Delete()
    }
}

# ---------------------------------------------------------------------
# 2) "... total embodiment of ethics for all of (STEM) with proven ..."
#    -> "... total embodiment of ethics, (STEM) and systems development
#         and engineering with proven ..."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "ethics for all of (STEM) with proven", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "ethics, (STEM) and systems development and engineering with proven", 2
) | Out-Null

# ---------------------------------------------------------------------
# 3) " Nuclear & Aerospace Engineer" (the one right after "CLM ~ Concepts
#    Learning Machine LLC." - it has a leading space) splits into three
#    runs: " Classic & Quantum Mechanics", " ", "Engineer".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    " Nuclear & Aerospace Engineer", $true, $false, $false, $false, $false,
    $true, 1, $false,
    " Classic & Quantum Mechanics Engineer", 1
) | Out-Null

$rng2 = $d.Content
$rng2.Find.Execute(
    "Mechanics Engineer", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
) | Out-Null
$spaceStart = $rng2.Start + 9
$spaceEnd = $spaceStart + 1
$spaceRng = $d.Range($spaceStart, $spaceEnd)
$spaceRng.Font.Bold = 0

# ---------------------------------------------------------------------
# 4) "Science and engineering with ethical algorithmic procedures."
#    -> "... ethical algorithmic processes & procedures."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "ethical algorithmic procedures.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "ethical algorithmic processes & procedures.", 2
) | Out-Null

# ---------------------------------------------------------------------
# 5) Reword the WaveLength & Impedance Meter sentence.
# ---------------------------------------------------------------------
$oldWave = "The Natural Human WaveLength & Impedance Meter/Visual Recognition Medical Instrument utility invention I engineered and built, submitted a white paper to the U.S Army Research Laboratory(ARL)."
$newWave = "The working utility invention Natural Human WaveLength & Impedance Meter/Visual Recognition Medical Instrument I engineered and built, submitted a white paper to the U.S Army Research Laboratory(ARL)."
$d.Content.Find.Execute($oldWave, $true, $false, $false, $false, $false, $true, 1, $false, $newWave, 2) | Out-Null

# ---------------------------------------------------------------------
# 6) "Engineered and built Toke Core; the hyperprotovisor software running"
#    -> "Engineered and built Toke Core the working hyperprotovisor systems
#         software running"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Toke Core; the hyperprotovisor software running", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Toke Core the working hyperprotovisor systems software running", 2
) | Out-Null

# ---------------------------------------------------------------------
# 7) "... building a Nuclear Fusion Reactor and a Fusion Reactionary
#    Engine, ..." -> "... building the working utility invention Nuclear
#    Fusion Reactor and a working Fusion Reactionary Engine, ..."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "building a Nuclear Fusion Reactor and a Fusion Reactionary Engine,", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "building the working utility invention Nuclear Fusion Reactor and a working Fusion Reactionary Engine,", 2
) | Out-Null

# ---------------------------------------------------------------------
# 8) "Engineered and partially built the Full-Spatial ..."
#    -> "Engineered and partially built the working utility invention
#         Full-Spatial ..."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "partially built the Full-Spatial", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "partially built the working utility invention Full-Spatial", 2
) | Out-Null

Write-Output "Done"
